$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) tweaks ---
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# --- Row 2 (CON) tweaks: drop B2:D2, update E2 ---
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -6.811486582493159

# --- Row 3 (STR) tweaks: drop B3, update C3, add D3, update E3 ---
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 1.1256232347222057
$ws.Range("D3").Value = -13.535893596395896
$ws.Range("E3").Value = 21.287391935224136

# --- Selection now only spans the edited block ---
$ws.Range("B1:E3").Select() | Out-Null
